# Auto-generated edit script to update the cryptos price table
# matching the target diff (price/volume refresh + two row swaps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain (non-numeric-looking) string - safe to set directly.
$safeValues = @{
    'D2' = '22.467.95'
    'E2' = '  +0.31%  '
    'D3' = '1.573.57'
    'E3' = '  +0.10%  '
    'E4' = '  -0.16%  '
    'E5' = '  -0.14%  '
    'E6' = '  +0.08%  '
    'E7' = '  -0.52%  '
    'E8' = '  +0.01%  '
    'E9' = '  -0.17%  '
    'B10' = 'Polygon'
    'C10' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E10' = '  -1.16%  '
    'B11' = 'Dogecoin'
    'C11' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'E11' = '  -1.06%  '
    'E12' = '  -0.14%  '
    'E14' = '  +0.58%  '
    'E15' = '  +0.61%  '
    'D16' = '1.572.18'
    'E16' = '  -0.06%  '
    'E17' = '  -0.96%  '
    'E18' = '  +0.71%  '
    'E19' = '  -0.48%  '
    'E20' = '  -0.22%  '
    'E21' = '  +1.23%  '
    'E22' = '  -1.80%  '
    'E23' = '  +1.63%  '
    'D24' = '22.466.43'
    'E24' = '  +0.33%  '
    'E25' = '  -5.09%  '
    'E26' = '  -5.16%  '
    'E27' = '  -0.63%  '
    'E28' = '  +2.15%  '
    'E29' = '  -1.10%  '
    'E30' = '  +0.11%  '
    'D31' = '1.747.42'
    'E31' = '  -0.12%  '
    'E32' = '  +3.47%  '
    'E33' = '  -0.49%  '
    'E34' = '  -1.55%  '
    'E35' = '  -1.69%  '
    'E36' = '  -1.32%  '
    'E37' = '  +4.05%  '
    'E38' = '  -3.54%  '
    'E39' = '  -0.74%  '
    'E40' = '  +0.30%  '
    'E41' = '  +1.36%  '
    'E42' = '  -1.29%  '
    'E43' = '  -2.56%  '
    'E44' = '  -0.28%  '
    'E45' = '  -0.73%  '
    'E46' = '  +0.67%  '
    'E47' = '  -2.52%  '
    'B48' = 'NEARProtocol'
    'C48' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E48' = '  +0.15%  '
    'B49' = 'Quant'
    'C49' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E49' = '  +3.69%  '
    'E50' = '  -5.01%  '
    'E51' = '  -0.03%  '
}

# Cells whose new value looks like a plain number (e.g. '291.24') but must stay
# stored as text, matching the original inlineStr cell type in the workbook.
# Force a text number format first so Excel doesn't silently convert them to numbers
# (which would also strip meaningful trailing zeros, e.g. '91.00' -> 91).
$textValues = @{
    'D6' = '291.24'
    'D7' = '0.3741'
    'D8' = '49.89'
    'D10' = '1.147'
    'D11' = '0.07575'
    'D14' = '5.993'
    'D15' = '6.950'
    'D18' = '91.00'
    'D19' = '0.06736'
    'D21' = '6.287'
    'D22' = '16.45'
    'D23' = '12.21'
    'D25' = '2.327'
    'D26' = '2.598'
    'D27' = '20.15'
    'D28' = '148.44'
    'D29' = '5.003'
    'D30' = '126.00'
    'D32' = '1.047'
    'D33' = '6.149'
    'D34' = '1.984'
    'D35' = '9.885'
    'D36' = '0.08432'
    'D37' = '1.387'
    'D39' = '0.2299'
    'D41' = '5.504'
    'D43' = '0.6289'
    'D44' = '1.002'
    'D45' = '13.96'
    'D46' = '3.814'
    'D47' = '0.5875'
    'D48' = '2.093'
    'D49' = '129.79'
    'D51' = '0.07330'
}

foreach ($cell in $safeValues.Keys) {
    $ws.Range($cell).Value = $safeValues[$cell]
}

foreach ($cell in $textValues.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = '@'
    $rng.Value = $textValues[$cell]
}

Write-Output "Applied all cell updates."
